$d = $word.ActiveDocument

# "Parameters" > "Required & Optional in a Text Code": update the document
# title (first heading, inside the blue banner table) from the old spec-style
# title to the new Construct-Drafts style title used across the converted
# MarkDown docs.
$d.Content.Find.Execute(
    "Circle Language Spec: Parameters", $true, $false, $false, $false, $false,
    $true, 1, $false, "Circle Language Construct Drafts | Parameters", 2
)
